# "player switching + camera"
#
# Adds a new logged time-tracking entry (row 8) for working on the
# player-swap ability and camera, plus the formatting touch-up (H4 gets
# the new accounting-style number format) that came along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New number-format style (numFmtId 4, "#,##0.00") shows up as cellXfs
# entry #6 after this; applied to the stray H4 cell exactly like the
# author's workbook.
$ws.Range("H4").NumberFormat = "#,##0.00"

# Row 8 was blank before; copy the previous (fully styled) row's
# formatting down first so every column picks up the right style
# (s="1" for text/time cells, s="3" for the date, s="4" for the
# duration), then overwrite with the real content for this entry.
$ws.Range("A7:F7").Copy($ws.Range("A8"))

$ws.Range("A8").Value = "Aris"
$ws.Range("B8").Value = 45316
$ws.Range("C8").Formula = "=19+42/60"
$ws.Range("D8").Formula = "=21 + 30/60"
$ws.Range("E8").Formula = "=D8-C8"
$ws.Range("F8").Value = "Player Swap ability +  Camera"

# Restore the cursor to where the author left it.
[void]$ws.Range("L17").Select()
